$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 431-432, shifting existing rows 431..532 down to 433..534.
$ws.Rows("431:432").Insert()

# Row 431 - new weekly record (Primera)
$ws.Cells.Item(431, 1).Value = 6
$ws.Cells.Item(431, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(431, 3).Value = "Metropolitana"
$ws.Cells.Item(431, 4).Value = 44543
$ws.Cells.Item(431, 5).Value = 13
$ws.Cells.Item(431, 6).Value = 100112008
$ws.Cells.Item(431, 7).Value = "Coliflor"
$ws.Cells.Item(431, 8).Value = "Sin especificar"
$ws.Cells.Item(431, 9).Value = "Primera"
$ws.Cells.Item(431, 10).Value = 6700
$ws.Cells.Item(431, 11).Value = 700
$ws.Cells.Item(431, 12).Value = 750
$ws.Cells.Item(431, 13).Value = 719
$ws.Cells.Item(431, 14).Value = "`$/unidad"
$ws.Cells.Item(431, 15).Value = "Región Metropolitana"
$ws.Cells.Item(431, 16).Value = 719
$ws.Cells.Item(431, 17).Value = 1
$ws.Cells.Item(431, 18).Value = "Hortaliza"

# Row 432 - new weekly record (Segunda)
$ws.Cells.Item(432, 1).Value = 6
$ws.Cells.Item(432, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(432, 3).Value = "Metropolitana"
$ws.Cells.Item(432, 4).Value = 44543
$ws.Cells.Item(432, 5).Value = 13
$ws.Cells.Item(432, 6).Value = 100112008
$ws.Cells.Item(432, 7).Value = "Coliflor"
$ws.Cells.Item(432, 8).Value = "Sin especificar"
$ws.Cells.Item(432, 9).Value = "Segunda"
$ws.Cells.Item(432, 10).Value = 2100
$ws.Cells.Item(432, 11).Value = 600
$ws.Cells.Item(432, 12).Value = 600
$ws.Cells.Item(432, 13).Value = 600
$ws.Cells.Item(432, 14).Value = "`$/unidad"
$ws.Cells.Item(432, 15).Value = "Región Metropolitana"
$ws.Cells.Item(432, 16).Value = 600
$ws.Cells.Item(432, 17).Value = 1
$ws.Cells.Item(432, 18).Value = "Hortaliza"
